$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Genero value for Felipe's row from "M" to lowercase "m"
$ws.Range("C2").Value = "m"

# Move selection to E5 (as in the final saved state)
$ws.Range("E5").Select()
